$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "CORREO" -- copy the header style (fill+border) used by the rest of row 1
$ws.Range("H1").Value2 = "CORREO"
$ws.Range("H1").Interior.Pattern = $ws.Range("G1").Interior.Pattern
$ws.Range("H1").Interior.Color = $ws.Range("G1").Interior.Color
$ws.Range("H1").Borders(7).LineStyle = 1
$ws.Range("H1").Borders(8).LineStyle = 1
$ws.Range("H1").Borders(9).LineStyle = 1
$ws.Range("H1").Borders(10).LineStyle = 1

# Hyperlink cell H2 with the mail address
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:lcisneros@sempiterno-group.com", "", "", "lcisneros@sempiterno-group.com")
$ws.Range("H2").Borders(7).LineStyle = 1
$ws.Range("H2").Borders(8).LineStyle = 1
$ws.Range("H2").Borders(9).LineStyle = 1
$ws.Range("H2").Borders(10).LineStyle = 1

# G2 no longer needs its right edge border (H2 supplies the shared line)
$ws.Range("G2").Borders(10).LineStyle = 0

# Column sizing to fit the new data
$ws.Columns("G").ColumnWidth = 14.85546875
$ws.Columns("H").ColumnWidth = 31.7109375

# Selection as left by the editing session
$ws.Range("H3").Select()
